$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.091330051422119
$ws.Range("B1").Value = 2.404870271682739
$ws.Range("C1").Value = 4.97639274597168
$ws.Range("D1").Value = 2.278491258621216
$ws.Range("E1").Value = 1.288783073425293
